{"js": "// Apply all text replacements described by the diff.\n// Each (oldText, newText) pair is unique within the document, so we\n// resolve every search hit against the ORIGINAL text first, then apply\n// the text replacements -- this avoids cascading edits where an earlier\n// replacement's output would accidentally match a later search pattern\n// (this matters for the 3-way label rotation).\nconst pairs = [\n  [\"N = 90,243\", \"N = 90,222\"],\n  [\"MVPA min/week - Activity count, Median (Q1, Q3)\", \"MVPA min/week - Machine learning, Median (Q1, Q3)\"],\n  [\"MVPA min/week - Self-report (IPAQ), Median (Q1, Q3)\", \"MVPA min/week - Activity count, Median (Q1, Q3)\"],\n  [\"MVPA min/week - Machine learning, Median (Q1, Q3)\", \"MVPA min/week - Self-report (IPAQ), Median (Q1, Q3)\"],\n  [\"2,017 (2.2)\", \"2,016 (2.2)\"],\n  [\"1,284 (1.4)\", \"1,283 (1.4)\"],\n  [\"2,680 (3.0)\", \"2,679 (3.0)\"],\n  [\"87,563 (97)\", \"87,543 (97)\"],\n  [\"51,728 (57)\", \"51,718 (57)\"],\n  [\"38,515 (43)\", \"38,504 (43)\"],\n  [\"7,116 (7.9)\", \"7,114 (7.9)\"],\n  [\"22,113 (25)\", \"22,110 (25)\"],\n  [\"21,317 (24)\", \"21,314 (24)\"],\n  [\"39,697 (44)\", \"39,684 (44)\"],\n  [\"11,583 (13)\", \"11,581 (13)\"],\n  [\"19,658 (22)\", \"19,656 (22)\"],\n  [\"23,666 (26)\", \"23,661 (26)\"],\n  [\"20,900 (23)\", \"20,893 (23)\"],\n  [\"6,073 (6.7)\", \"6,071 (6.7)\"],\n  [\"8,363 (9.3)\", \"8,360 (9.3)\"],\n  [\"52,100 (58)\", \"52,089 (58)\"],\n  [\"32,066 (36)\", \"32,056 (36)\"],\n  [\"4,961 (5.5)\", \"4,959 (5.5)\"],\n  [\"18,240 (20)\", \"18,236 (20)\"],\n  [\"22,702 (25)\", \"22,695 (25)\"],\n  [\"23,631 (26)\", \"23,626 (26)\"],\n  [\"20,709 (23)\", \"20,706 (23)\"],\n  [\"64,272 (71)\", \"64,254 (71)\"],\n  [\"22,776 (25)\", \"22,773 (25)\"],\n  [\"16,294 (18)\", \"16,289 (18)\"],\n  [\"30,600 (34)\", \"30,593 (34)\"],\n  [\"43,349 (48)\", \"43,340 (48)\"],\n  [\"74,534 (83)\", \"74,514 (83)\"],\n  [\"14,988 (17)\", \"14,987 (17)\"],\n  [\"76,538 (85)\", \"76,521 (85)\"],\n  [\"13,281 (15)\", \"13,277 (15)\"],\n  [\"29,458 (33)\", \"29,453 (33)\"],\n  [\"54,043 (60)\", \"54,031 (60)\"],\n  [\"6,742 (7.5)\", \"6,738 (7.5)\"],\n  [\"15,369 (17)\", \"15,368 (17)\"],\n  [\"66,466 (74)\", \"66,451 (74)\"],\n  [\"8,408 (9.3)\", \"8,403 (9.3)\"],\n];\n\nconst body = context.document.body;\nconst pending = [];\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  pending.push({ oldText, newText, results });\n}\nawait context.sync();\n\nfor (const { oldText, newText, results } of pending) {\n  if (results.items.length !== 1) {\n    throw new Error(`expected exactly 1 match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n}\nawait context.sync();", "ps1": "# Apply all text replacements described by the diff.\n# Each (oldText, newText) pair is unique within the document. We resolve\n# every match's ABSOLUTE [start,end) offsets against the ORIGINAL content\n# first (Find.Execute scans forward in document order, including inside\n# table cells), and only then rewrite the ranges -- processed from the\n# LAST match to the FIRST -- so that an earlier edit's length delta never\n# invalidates the offsets of a not-yet-applied later edit. This also keeps\n# the 3-way label rotation (Activity count -> Machine learning ->\n# Self-report (IPAQ) -> Activity count) from cascading: offsets are fixed\n# up front, before any text is rewritten.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"N = 90,243\", \"N = 90,222\")\n    ,@(\"MVPA min/week - Activity count, Median (Q1, Q3)\", \"MVPA min/week - Machine learning, Median (Q1, Q3)\")\n    ,@(\"MVPA min/week - Self-report (IPAQ), Median (Q1, Q3)\", \"MVPA min/week - Activity count, Median (Q1, Q3)\")\n    ,@(\"MVPA min/week - Machine learning, Median (Q1, Q3)\", \"MVPA min/week - Self-report (IPAQ), Median (Q1, Q3)\")\n    ,@(\"2,017 (2.2)\", \"2,016 (2.2)\")\n    ,@(\"1,284 (1.4)\", \"1,283 (1.4)\")\n    ,@(\"2,680 (3.0)\", \"2,679 (3.0)\")\n    ,@(\"87,563 (97)\", \"87,543 (97)\")\n    ,@(\"51,728 (57)\", \"51,718 (57)\")\n    ,@(\"38,515 (43)\", \"38,504 (43)\")\n    ,@(\"7,116 (7.9)\", \"7,114 (7.9)\")\n    ,@(\"22,113 (25)\", \"22,110 (25)\")\n    ,@(\"21,317 (24)\", \"21,314 (24)\")\n    ,@(\"39,697 (44)\", \"39,684 (44)\")\n    ,@(\"11,583 (13)\", \"11,581 (13)\")\n    ,@(\"19,658 (22)\", \"19,656 (22)\")\n    ,@(\"23,666 (26)\", \"23,661 (26)\")\n    ,@(\"20,900 (23)\", \"20,893 (23)\")\n    ,@(\"6,073 (6.7)\", \"6,071 (6.7)\")\n    ,@(\"8,363 (9.3)\", \"8,360 (9.3)\")\n    ,@(\"52,100 (58)\", \"52,089 (58)\")\n    ,@(\"32,066 (36)\", \"32,056 (36)\")\n    ,@(\"4,961 (5.5)\", \"4,959 (5.5)\")\n    ,@(\"18,240 (20)\", \"18,236 (20)\")\n    ,@(\"22,702 (25)\", \"22,695 (25)\")\n    ,@(\"23,631 (26)\", \"23,626 (26)\")\n    ,@(\"20,709 (23)\", \"20,706 (23)\")\n    ,@(\"64,272 (71)\", \"64,254 (71)\")\n    ,@(\"22,776 (25)\", \"22,773 (25)\")\n    ,@(\"16,294 (18)\", \"16,289 (18)\")\n    ,@(\"30,600 (34)\", \"30,593 (34)\")\n    ,@(\"43,349 (48)\", \"43,340 (48)\")\n    ,@(\"74,534 (83)\", \"74,514 (83)\")\n    ,@(\"14,988 (17)\", \"14,987 (17)\")\n    ,@(\"76,538 (85)\", \"76,521 (85)\")\n    ,@(\"13,281 (15)\", \"13,277 (15)\")\n    ,@(\"29,458 (33)\", \"29,453 (33)\")\n    ,@(\"54,043 (60)\", \"54,031 (60)\")\n    ,@(\"6,742 (7.5)\", \"6,738 (7.5)\")\n    ,@(\"15,369 (17)\", \"15,368 (17)\")\n    ,@(\"66,466 (74)\", \"66,451 (74)\")\n    ,@(\"8,408 (9.3)\", \"8,403 (9.3)\")\n)\n\n$resolved = @()\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $searchRange = $d.Content\n    $found = $searchRange.Find.Execute($oldText, $false, $true)\n    if (-not $found) {\n        throw \"Could not find expected text: $oldText\"\n    }\n    $resolved += , @($searchRange.Start, $searchRange.End, $oldText, $newText)\n}\n\nfor ($i = $resolved.Count - 1; $i -ge 0; $i--) {\n    $item = $resolved[$i]\n    $s = $item[0]\n    $e = $item[1]\n    $oldText = $item[2]\n    $newText = $item[3]\n    $rng = $d.Range($s, $e)\n    if ($rng.Text -ne $oldText) {\n        throw \"Range text mismatch at [$s,$e): expected [$oldText] got [$($rng.Text)]\"\n    }\n    $rng.Text = $newText\n}"}
